$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.230.66"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").Value = "2.482.33"
$ws.Range("E3").Value = "  +1.02%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.78"
$ws.Range("E5").Value = "  +1.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.04"
$ws.Range("E6").Value = "  +3.61%  "
$ws.Range("E8").Value = "  +0.62%  "
$ws.Range("D9").Value = "2.482.43"
$ws.Range("E9").Value = "  +1.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.139"
$ws.Range("E10").Value = "  +3.75%  "
$ws.Range("E11").Value = "  +1.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.94"
$ws.Range("E12").Value = "  +1.13%  "
$ws.Range("E13").Value = "  +0.66%  "
$ws.Range("D14").Value = "2.936.69"
$ws.Range("E14").Value = "  +1.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.57"
$ws.Range("E15").Value = "  +1.13%  "
$ws.Range("D16").Value = "67.136.22"
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("D18").Value = "2.499.20"
$ws.Range("E18").Value = "  +1.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.67"
$ws.Range("E19").Value = "  -0.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.02"
$ws.Range("E20").Value = "  -2.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "352.20"
$ws.Range("E21").Value = "  -0.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.04"
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.06"
$ws.Range("E24").Value = "  -0.38%  "
$ws.Range("E25").Value = "  +1.01%  "
$ws.Range("E26").Value = "  +3.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.25"
$ws.Range("E27").Value = "  +4.42%  "
$ws.Range("D28").Value = "2.600.48"
$ws.Range("E28").Value = "  +0.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("D30").Value = "0.0₃0914"
$ws.Range("E30").Value = "  +2.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "509.57"
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.73"
$ws.Range("E32").Value = "  -0.66%  "
$ws.Range("E33").Value = "  +2.91%  "
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "162.44"
$ws.Range("E36").Value = "  +2.17%  "
$ws.Range("E37").Value = "  +2.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.73"
$ws.Range("E38").Value = "  +0.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.16"
$ws.Range("E39").Value = "  -1.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.33"
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("E42").Value = "  +2.30%  "
$ws.Range("E43").Value = "  +1.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.83"
$ws.Range("E44").Value = "  +1.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.38"
$ws.Range("E45").Value = "  +3.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "143.99"
$ws.Range("E46").Value = "  +2.02%  "
$ws.Range("D47").Value = "0.0₆0261"
$ws.Range("E47").Value = "  +4.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.49"
$ws.Range("E48").Value = "  +0.62%  "
$ws.Range("E49").Value = "  +0.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0735"
$ws.Range("E50").Value = "  +0.41%  "
$ws.Range("E51").Value = "  -0.13%  "
